$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8151398301124573
$ws.Range("B1").Value = 1.640724897384644
$ws.Range("C1").Value = 2.834611654281616
$ws.Range("D1").Value = 3.610299348831177
$ws.Range("E1").Value = 2.286367177963257
